{"js": "// The document contains a single 20-row x 5-column table of simple\n// arithmetic problems (e.g. \"29+11=\"). The edit replaces the text of every\n// one of the 100 cells with a new problem, keeping the table's shape\n// (20 rows x 5 columns) and all existing formatting untouched.\nconst newValues = [\n  [\"66-12=\", \"34-26=\", \"76-49=\", \"22+42=\", \"87-17=\"],\n  [\"27+58=\", \"93-52=\", \"87-13=\", \"80-1=\", \"80-67=\"],\n  [\"23+65=\", \"79-64=\", \"61-0=\", \"74-7=\", \"53+45=\"],\n  [\"1+25=\", \"62+15=\", \"59+34=\", \"83-25=\", \"73-11=\"],\n  [\"83-54=\", \"70-6=\", \"46-3=\", \"81+8=\", \"47-3=\"],\n  [\"39-24=\", \"84-39=\", \"86+13=\", \"99-18=\", \"38-32=\"],\n  [\"24+60=\", \"65-55=\", \"14+56=\", \"58-11=\", \"4+9=\"],\n  [\"74-17=\", \"54+20=\", \"94+4=\", \"53-37=\", \"55-33=\"],\n  [\"80-27=\", \"93-83=\", \"80-2=\", \"66+6=\", \"1+27=\"],\n  [\"85-77=\", \"38-10=\", \"23+58=\", \"88-53=\", \"99-60=\"],\n  [\"88-84=\", \"25+29=\", \"34-23=\", \"83-50=\", \"82+16=\"],\n  [\"59-14=\", \"96-22=\", \"46+30=\", \"90-82=\", \"36+57=\"],\n  [\"2+42=\", \"74+14=\", \"6+23=\", \"35+0=\", \"53-19=\"],\n  [\"50+41=\", \"98-40=\", \"62+2=\", \"10+9=\", \"76-37=\"],\n  [\"74-16=\", \"13+36=\", \"85-6=\", \"43+6=\", \"89-44=\"],\n  [\"31-11=\", \"69-24=\", \"50-46=\", \"18+7=\", \"3+45=\"],\n  [\"9+87=\", \"66-53=\", \"47+23=\", \"75+14=\", \"76-15=\"],\n  [\"6+70=\", \"60-41=\", \"66-29=\", \"92-35=\", \"91-88=\"],\n  [\"45+36=\", \"90-77=\", \"72+1=\", \"56+35=\", \"73+22=\"],\n  [\"56+18=\", \"42-39=\", \"84-0=\", \"79-16=\", \"31+14=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = Math.min(table.rowCount, newValues.length);\nfor (let r = 0; r < rowCount; r++) {\n  const rowValues = newValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = rowValues[c];\n  }\n}\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of simple\n# arithmetic problems (e.g. \"29+11=\"). The edit replaces the text of every\n# one of the 100 cells with a new problem, keeping the table's shape\n# (20 rows x 5 columns) and all existing formatting untouched.\n$newValues = @(\n  @(\"66-12=\", \"34-26=\", \"76-49=\", \"22+42=\", \"87-17=\"),\n  @(\"27+58=\", \"93-52=\", \"87-13=\", \"80-1=\", \"80-67=\"),\n  @(\"23+65=\", \"79-64=\", \"61-0=\", \"74-7=\", \"53+45=\"),\n  @(\"1+25=\", \"62+15=\", \"59+34=\", \"83-25=\", \"73-11=\"),\n  @(\"83-54=\", \"70-6=\", \"46-3=\", \"81+8=\", \"47-3=\"),\n  @(\"39-24=\", \"84-39=\", \"86+13=\", \"99-18=\", \"38-32=\"),\n  @(\"24+60=\", \"65-55=\", \"14+56=\", \"58-11=\", \"4+9=\"),\n  @(\"74-17=\", \"54+20=\", \"94+4=\", \"53-37=\", \"55-33=\"),\n  @(\"80-27=\", \"93-83=\", \"80-2=\", \"66+6=\", \"1+27=\"),\n  @(\"85-77=\", \"38-10=\", \"23+58=\", \"88-53=\", \"99-60=\"),\n  @(\"88-84=\", \"25+29=\", \"34-23=\", \"83-50=\", \"82+16=\"),\n  @(\"59-14=\", \"96-22=\", \"46+30=\", \"90-82=\", \"36+57=\"),\n  @(\"2+42=\", \"74+14=\", \"6+23=\", \"35+0=\", \"53-19=\"),\n  @(\"50+41=\", \"98-40=\", \"62+2=\", \"10+9=\", \"76-37=\"),\n  @(\"74-16=\", \"13+36=\", \"85-6=\", \"43+6=\", \"89-44=\"),\n  @(\"31-11=\", \"69-24=\", \"50-46=\", \"18+7=\", \"3+45=\"),\n  @(\"9+87=\", \"66-53=\", \"47+23=\", \"75+14=\", \"76-15=\"),\n  @(\"6+70=\", \"60-41=\", \"66-29=\", \"92-35=\", \"91-88=\"),\n  @(\"45+36=\", \"90-77=\", \"72+1=\", \"56+35=\", \"73+22=\"),\n  @(\"56+18=\", \"42-39=\", \"84-0=\", \"79-16=\", \"31+14=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = [Math]::Min($t.Rows.Count, $newValues.Count)\nfor ($r = 1; $r -le $rowCount; $r++) {\n  $rowValues = $newValues[$r - 1]\n  for ($c = 1; $c -le $rowValues.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $rowValues[$c - 1]\n  }\n}\n"}
